$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2319.7144
$ws.Range("J17").Value = 2319.7144
$ws.Range("L17").Value = 6959.1432
$ws.Range("N17").Value = -7295.1432
$ws.Range("H40").Value = 1891.6666
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 1855.5555
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 1855.5555
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2205.5555
$ws.Range("H64").Value = 3411.9412
$ws.Range("I64").Value = 3328.5715
$ws.Range("K64").Value = 3328.5715
$ws.Range("M64").Value = -3080.5715
$ws.Range("H67").Value = 3411.9412
$ws.Range("I67").Value = 3328.5715
$ws.Range("K67").Value = 3328.5715
$ws.Range("M67").Value = -2470.5715
$ws.Range("H74").Value = 3591.8147
$ws.Range("I74").Value = 3031
$ws.Range("J74").Value = 3872.2222
$ws.Range("K74").Value = 3031
$ws.Range("L74").Value = 3872.2222
$ws.Range("M74").Value = -2095
$ws.Range("N74").Value = -5744.2222
$ws.Range("H77").Value = 3591.8147
$ws.Range("I77").Value = 3031
$ws.Range("J77").Value = 3872.2222
$ws.Range("K77").Value = 15155
$ws.Range("L77").Value = 19361.111
$ws.Range("M77").Value = -10475
$ws.Range("N77").Value = -28721.111
$ws.Range("H96").Value = 1002.5217
$ws.Range("I96").Value = 1044.5
$ws.Range("J96").Value = 956.7273
$ws.Range("K96").Value = 3133.5
$ws.Range("L96").Value = 2870.1819
$ws.Range("M96").Value = -1760.5
$ws.Range("N96").Value = -5616.1819
$ws.Range("H113").Value = 2433.3333
$ws.Range("I113").Value = 2464.2856
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2464.2856
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 789.7143999999998
$ws.Range("N113").Value = -8508
$ws.Range("H135").Value = 1529.3043
$ws.Range("I135").Value = 1260.5161
$ws.Range("J135").Value = 2084.8
$ws.Range("K135").Value = 11344.6449
$ws.Range("L135").Value = 18763.2
$ws.Range("M135").Value = -8809.644900000001
$ws.Range("N135").Value = -23833.2
$ws.Range("H138").Value = 2206.074
$ws.Range("I138").Value = 1374.2142
$ws.Range("J138").Value = 3101.923
$ws.Range("K138").Value = 4122.642599999999
$ws.Range("L138").Value = 9305.769
$ws.Range("M138").Value = 1017.357400000001
$ws.Range("N138").Value = -19585.769
$ws.Range("H139").Value = 52354
$ws.Range("J139").Value = 52354
$ws.Range("L139").Value = 52354
$ws.Range("N139").Value = -62634
$ws.Range("H140").Value = 49933.332
$ws.Range("J140").Value = 49933.332
$ws.Range("L140").Value = 49933.332
$ws.Range("N140").Value = -60293.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8655.839
$ws.Range("I32").Value = 8506.17
$ws.Range("K32").Value = 8506.17
$ws.Range("M32").Value = -8219.17
$ws.Range("H45").Value = 1710.7646
$ws.Range("I45").Value = 1794
$ws.Range("J45").Value = 1511
$ws.Range("K45").Value = 1794
$ws.Range("L45").Value = 1511
$ws.Range("M45").Value = -1417
$ws.Range("N45").Value = -2265
$ws.Range("H61").Value = 18520176
$ws.Range("I61").Value = 22729016
$ws.Range("J61").Value = 1272.8
$ws.Range("K61").Value = 22729016
$ws.Range("L61").Value = 1272.8
$ws.Range("M61").Value = -22728804
$ws.Range("N61").Value = -1696.8
$ws.Range("H74").Value = 14709112
$ws.Range("I74").Value = 25001922
$ws.Range("J74").Value = 5099
$ws.Range("K74").Value = 25001922
$ws.Range("L74").Value = 5099
$ws.Range("M74").Value = -25001048
$ws.Range("N74").Value = -6847
$ws.Range("H77").Value = 14709112
$ws.Range("I77").Value = 25001922
$ws.Range("J77").Value = 5099
$ws.Range("K77").Value = 125009610
$ws.Range("L77").Value = 25495
$ws.Range("M77").Value = -125005242
$ws.Range("N77").Value = -34231
$ws.Range("H136").Value = 18520176
$ws.Range("I136").Value = 22729016
$ws.Range("J136").Value = 1272.8
$ws.Range("K136").Value = 68187048
$ws.Range("L136").Value = 3818.4
$ws.Range("M136").Value = -68184498
$ws.Range("N136").Value = -8918.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1220.1428
$ws.Range("I16").Value = 1078.2
$ws.Range("J16").Value = 1575
$ws.Range("K16").Value = 1078.2
$ws.Range("L16").Value = 1575
$ws.Range("M16").Value = -791.2
$ws.Range("N16").Value = -2149
$ws.Range("H113").Value = 1220.1428
$ws.Range("I113").Value = 1078.2
$ws.Range("J113").Value = 1575
$ws.Range("K113").Value = 1078.2
$ws.Range("L113").Value = 1575
$ws.Range("M113").Value = 1091.8
$ws.Range("N113").Value = -5915
$ws.Range("H122").Value = 1696.1305
$ws.Range("I122").Value = 1791.1666
$ws.Range("J122").Value = 1354
$ws.Range("K122").Value = 5373.4998
$ws.Range("L122").Value = 4062
$ws.Range("M122").Value = -2923.4998
$ws.Range("N122").Value = -8962
$ws.Range("H134").Value = 1797.8422
$ws.Range("I134").Value = 1797.8422
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5393.5266
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2858.5266
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 76300
$ws.Range("J135").Value = 76300
$ws.Range("L135").Value = 76300
$ws.Range("N135").Value = -86440
$ws.Range("H140").Value = 37325.92
$ws.Range("J140").Value = 37325.92
$ws.Range("L140").Value = 37325.92
$ws.Range("N140").Value = -47685.92
$ws.Range("H141").Value = 106883.695
$ws.Range("J141").Value = 106883.695
$ws.Range("L141").Value = 106883.695
$ws.Range("N141").Value = -117243.695

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 4300
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16352
$ws.Range("H100").Value = 8715
$ws.Range("J100").Value = 8715
$ws.Range("L100").Value = 26145
$ws.Range("N100").Value = -27767
$ws.Range("H103").Value = 5491.6665
$ws.Range("I103").Value = 627.25
$ws.Range("J103").Value = 6636.2354
$ws.Range("K103").Value = 1881.75
$ws.Range("L103").Value = 19908.7062
$ws.Range("M103").Value = -1002.75
$ws.Range("N103").Value = -21666.7062
$ws.Range("H106").Value = 9512.857
$ws.Range("J106").Value = 9512.857
$ws.Range("L106").Value = 28538.571
$ws.Range("N106").Value = -30430.571
$ws.Range("H112").Value = 8760.125
$ws.Range("I112").Value = 7360.3335
$ws.Range("J112").Value = 9600
$ws.Range("K112").Value = 22081.0005
$ws.Range("L112").Value = 28800
$ws.Range("M112").Value = -20973.0005
$ws.Range("N112").Value = -31016
$ws.Range("H134").Value = 3672.2222
$ws.Range("I134").Value = 1640.625
$ws.Range("J134").Value = 6627.273
$ws.Range("K134").Value = 4921.875
$ws.Range("L134").Value = 19881.819
$ws.Range("M134").Value = 148.125
$ws.Range("N134").Value = -30021.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2519.3333
$ws.Range("I97").Value = 2467.7778
$ws.Range("J97").Value = 2674
$ws.Range("K97").Value = 2467.7778
$ws.Range("L97").Value = 2674
$ws.Range("M97").Value = -1971.7778
$ws.Range("N97").Value = -3666
$ws.Range("H122").Value = 13335994
$ws.Range("I122").Value = 22223390
$ws.Range("J122").Value = 4900
$ws.Range("K122").Value = 66670170
$ws.Range("L122").Value = 14700
$ws.Range("M122").Value = -66667720
$ws.Range("N122").Value = -19600
$ws.Range("H133").Value = 67857.14
$ws.Range("J133").Value = 67857.14
$ws.Range("L133").Value = 67857.14
$ws.Range("N133").Value = -77977.14
$ws.Range("H138").Value = 59249.25
$ws.Range("J138").Value = 59249.25
$ws.Range("L138").Value = 59249.25
$ws.Range("N138").Value = -69529.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6949722
$ws.Range("I132").Value = 3113.1277
$ws.Range("J132").Value = 20009346
$ws.Range("K132").Value = 9339.383099999999
$ws.Range("L132").Value = 60028038
$ws.Range("M132").Value = -6809.383099999999
$ws.Range("N132").Value = -60033098
$ws.Range("H139").Value = 58893
$ws.Range("J139").Value = 58893
$ws.Range("L139").Value = 58893
$ws.Range("N139").Value = -69173
